$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix E6: store "10000" as text instead of a number ---
$ws.Range("E6").Value = "'10000"

# --- Row 7 ---
$ws.Range("A7").Value = "'2345678990112444"
$ws.Range("B7").Value = "BG6701HI"
$ws.Range("C7").Value = "Tiara"
$ws.Range("D7").Value = "01-08-2025 06:23"
$ws.Range("E7").Value = "'0"
$ws.Range("F7").Value = "Bank Rakyat Indonesia (BRI)"

# --- Row 8 ---
$ws.Range("A8").Value = "'1234456278949542"
$ws.Range("B8").Value = "BG4576HI"
$ws.Range("C8").Value = "Nia"
$ws.Range("D8").Value = "01-08-2025 06:28"
$ws.Range("E8").Value = "'10000"
$ws.Range("F8").Value = "Bank Syariah Indonesia (BSI)"

# --- Row 9 ---
$ws.Range("A9").Value = "'1234456278949542"
$ws.Range("B9").Value = "BG4576HI"
$ws.Range("C9").Value = "Nia"
$ws.Range("D9").Value = "01-08-2025 06:28"
$ws.Range("E9").Value = 20000
$ws.Range("F9").Value = "Bank Mandiri"
